# Updated remaining queries for C3DC
# The SQL queries embedded in column B (and C2) joined df_participant /
# df_study / df_diagnoses / df_treatments / df_treatment_resp / df_survival /
# df_reference_files using the generic ".id" column names. The schema was
# updated to use explicit "study_id" / "participant_id" column names, so the
# JOIN predicates (and the quoted dotted aliases used by the duckdb/pandas
# style lookups) need to be updated to match everywhere they occur.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query([string]$cellAddr) {
    $cell = $ws.Range($cellAddr)
    $text = $cell.Value()
    if ($text -eq $null) { return }

    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

    $cell.Value = $text
}

# StatQuery (C2), StudiesTab query (B2), ParticipantsTab (B3), DiagnosisTab (B4),
# TreatmentTab (B5), TreatmentRespTab (B6), SurvivalTab (B7)
Update-Query "C2"
Update-Query "B2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# Column C was widened (and no longer relies on the old "best fit" measurement)
# to comfortably fit the updated, slightly longer query text.
$ws.Columns.Item(3).ColumnWidth = 69

# The author's cursor ended up back on B2 (top of the sheet) instead of C7.
$ws.Range("B2").Select()
